$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The energy buy-price table (column B) for rows 182-277 was re-linearized
# from a flat 0.2 to 0.24 as part of the MPC boiler model tuning.
$ws.Range("B182:B277").Value = 0.24

# Reflect the scrolled/selected view position captured for this edit.
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("D277").Select()
$excel.ActiveWindow.ScrollRow = 251
$excel.ActiveWindow.ScrollColumn = 1
